$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.47%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.254"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.06%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07659"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.64%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.622"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.84%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9193"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.00%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.465"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.91%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1246"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "11.63%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1836"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.12%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09091"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.23%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04360"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.43%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.13%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001262"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.89%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005783"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.68%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.007498"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,391.63%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.355"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.13%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.332"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.87%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.86%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.179"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.16%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1382"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.35%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2924"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.97%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04061"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.25%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004168"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.31%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.20%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02449"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2.61%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05296"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.21%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007846"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.89%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1314"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.51%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006817"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.02%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001903"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.48%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008351"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.17%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3337"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.56%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006903"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.17%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2055"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,897.95%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"

Write-Host "Applied all changes"
